$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("정동연")
$ws.Activate()

# --- Row 34: edit end-time / interrupt / delta values (brand/activity unchanged) ---
$ws.Range("C34").Value = 0.80555555555555547
$ws.Range("D34").Value = 170
$ws.Range("E34").Value = 480
$ws.Range("E34").NumberFormat = "0"

# --- Row 35: now reflects what used to be logged across rows 35 & 36 ---
$ws.Range("A35").Value = 43802
$ws.Range("B35").Value = 0.70138888888888884
$ws.Range("C35").Value = 0.79166666666666663
$ws.Range("D35").Value = 30
$ws.Range("E35").Value = 130

# --- Row 36: cleared out (kept formatting, no data) ---
$ws.Range("A36:F36").ClearContents()

# --- Rows 37/38: empty "date" cells now styled to match the data rows above ---
$ws.Range("A37").NumberFormat = 'm"월"\ d"일";@'
$ws.Range("A38").NumberFormat = 'm"월"\ d"일";@'

# --- sheet view: scroll back to top-left, move selection to F36 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("F36").Select()
